$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = 60300
    3 = 131276
    4 = 122717
    5 = 13295
    6 = 27133
    8 = 3021179
    9 = 217008
    10 = 29850
    11 = 623201
    12 = 321798
    13 = 10576
    14 = 180462
    15 = 763682
    16 = 3881
    17 = 361897
    18 = 996896
    19 = 12674
    20 = 7884
    21 = 1119
    22 = 308984
    24 = 47851
    25 = 14779529
    27 = 405825
    28 = 13324
    29 = 142842
    30 = 4046
    31 = 24548
    32 = 15361
    34 = 1251811
    37 = 1215815
    38 = 102549
    39 = 2905254
    40 = 3845
    41 = 11016
    42 = 30004
    43 = 256676
    44 = 46154
    45 = 335522
    46 = 109625
    47 = 66911
    48 = 1634619
    49 = 254368
    50 = 11201
    52 = 267681
    53 = 388046
    54 = 230713
    55 = 69727
    57 = 3673
    58 = 122943
    59 = 18460
    60 = 258813
    61 = 121
    62 = 87345
    63 = 5717160
    64 = 23201
    66 = 312954
    67 = 3438186
    68 = 92683
    69 = 348568
    71 = 228871
    72 = 22333
    74 = 13564
    76 = 214475
    77 = 784111
    78 = 6483
    79 = 20282833
    80 = 1682004
    81 = 2555587
    82 = 1079998
    83 = 250290
    84 = 838621
    85 = 4050708
    86 = 46039
    87 = 609625
    88 = 714173
    89 = 383164
    90 = 160559
    91 = 124269
    93 = 277832
    94 = 96337
    95 = 966
    96 = 119953
    97 = 528457
    99 = 2110
    100 = 178335
    101 = 2946
    102 = 250337
    103 = 67495
    104 = 37751
    105 = 34099
    106 = 417512
    107 = 31330
    108 = 13937
    109 = 30354
    111 = 18476
    113 = 2349900
    114 = 251378
    115 = 2465
    116 = 40396
    117 = 97718
    118 = 511912
    119 = 70000
    121 = 343418
    122 = 1543494
    123 = 2623
    125 = 5261
    126 = 165199
    127 = 113952
    128 = 196900
    129 = 837523
    130 = 365299
    131 = 11262
    132 = 284453
    133 = 1814127
    134 = 1062225
    135 = 2805756
    136 = 837457
    137 = 207592
    138 = 1058337
    139 = 4776844
    140 = 25311
    145 = 5067
    147 = 420301
    148 = 40464
    149 = 693169
    150 = 6373
    152 = 61235
    153 = 383228
    154 = 241883
    157 = 1584961
    158 = 10602
    159 = 3540430
    160 = 113676
    162 = 10543
    164 = 663952
    165 = 22977
    166 = 1145
    169 = 71025
    170 = 2524
    172 = 11471
    173 = 312747
    174 = 4900121
    175 = 32472201
    176 = 41973
    177 = 2140838
    178 = 525567
    179 = 4437505
    180 = 204120
    181 = 92006
    183 = 200931
    184 = 2985
    185 = 298203
    186 = 6363
    187 = 91722
    188 = 38293
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $updates[$row]
}

